$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update a few previously-missing / newly-missing values in column E
#    (header "D") for rows 19, 21 and 23.
$ws.Range("E19").Value = -6.5
$ws.Range("E21").ClearContents()
$ws.Range("E23").Value = -7

# 2) Drop the "RM 232" (row 26) and "SC 92" (row 28) data rows entirely,
#    shifting everything below them up. Delete the higher-numbered row
#    first so row 26's position is still valid when it is deleted next.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# 3) After the two deletions, "SC 101" has moved up to row 27 and
#    "SC 232" (the last row) has moved up to row 33. Update their
#    column E values to match the new data.
$ws.Range("E27").ClearContents()
$ws.Range("E33").Value = -10.7
